# ============================================================
# edit.ps1 - applies the 2022-Q1 data addition described by the diff:
#   1. Insert a new worksheet '2022-Q1' (built from the '2021-Q4'
#      template so header/styles match) directly before '总计'.
#   2. Fill it with the 34 rows of per-fund holding data for 2022-Q1.
#   3. Insert a new first data row into '总计' for 2022-Q1 and
#      bump the existing row indices by one.
# ============================================================

$wb = $excel.ActiveWorkbook

# ---- Step 1: create the new "2022-Q1" sheet right before "总计" ----
# Copying "2021-Q4" gives us the exact same header row / column styles
# (B..H headers, bold-centered-bordered style on column A) for free.
$totalSheet = $wb.Worksheets.Item("总计")
$templateSheet = $wb.Worksheets.Item("2021-Q4")
$templateSheet.Copy($totalSheet, $null)

# re-fetch by name: inserting the copy shifts "总计" one slot to the
# right, and the pre-copy $totalSheet reference keeps its stale .Index
$totalSheet = $wb.Worksheets.Item("总计")
$newSheet = $wb.Worksheets.Item($totalSheet.Index - 1)
$newSheet.Name = "2022-Q1"

# ---- Step 2: extend the copied sheet from 25 to 34 data rows ----
# Rows 2..26 already exist (copied from "2021-Q4"); rows 27..35 are new
# and need column A to inherit the same bold/border/center style used
# by the rest of column A, so clone it from row 26 first.
$newSheet.Range("A26").Copy($newSheet.Range("A27:A35"))

# Force columns B, D, E, F, G to be stored as literal text (matching the
# source data, which keeps fund codes / percentages / NAV figures as
# strings rather than coercing them to numbers).
$newSheet.Range("B2:B35").NumberFormat = "@"
$newSheet.Range("D2:G35").NumberFormat = "@"

$fundData = @(
    ,@('516970', '广发中证基建工程交易型开放式指数证券投资基金', '59.01', '99.38', '5.74', '3.3872', 5)
    ,@('161810', '银华内需精选混合(LOF)', '25.59', '94.71', '5.08', '1.3000', 10)
    ,@('009394', '银华同力精选混合', '20.03', '94.68', '5.66', '1.1337', 6)
    ,@('010963', '信达澳银周期动力混合', '30.52', '89.82', '3.11', '0.9492', 8)
    ,@('165525', '信诚中证基建工程指数（LOF）', '17.06', '94.00', '5.43', '0.9264', 5)
    ,@('970016', '中信建投价值增长混合A', '20.46', '78.12', '2.97', '0.6077', 8)
    ,@('011498', '富国沪深300基本面精选股票型证券投资基金A', '11.93', '90.96', '3.96', '0.4724', 7)
    ,@('516950', '银华中证基建交易型开放式指数证券投资基金', '10.41', '97.55', '4.26', '0.4435', 6)
    ,@('002780', '新疆前海联合泓鑫灵活配置混合A', '8.50', '75.30', '3.45', '0.2932', 7)
    ,@('005671', '新疆前海联合研究优选灵活配置混合A', '5.62', '77.46', '3.65', '0.2051', 8)
    ,@('001403', '招商国企改革主题混合', '3.57', '87.32', '5.11', '0.1824', 3)
    ,@('002504', '鹏华金鼎灵活配置混合A', '2.49', '77.53', '6.72', '0.1673', 3)
    ,@('180020', '银华成长先锋混合', '3.05', '79.81', '5.27', '0.1607', 7)
    ,@('000029', '富国宏观策略灵活配置混合', '5.97', '90.27', '2.58', '0.1540', 5)
    ,@('519770', '交银优择回报灵活配置混合A', '13.78', '21.65', '0.86', '0.1185', 4)
    ,@('970017', '中信建投价值增长混合C', '3.48', '78.12', '2.97', '0.1034', 8)
    ,@('160421', '华安智增精选灵活配置混合（LOF）', '2.11', '84.67', '4.22', '0.0890', 7)
    ,@('009907', '湘财长泽灵活配置混合A', '1.79', '80.11', '4.94', '0.0884', 3)
    ,@('011160', '富国质量成长6个月持有期混合A', '3.80', '91.55', '2.12', '0.0806', 10)
    ,@('007043', '新疆前海联合泓鑫灵活配置混合C', '1.37', '75.30', '3.45', '0.0473', 7)
    ,@('519771', '交银优择回报灵活配置混合C', '5.19', '21.65', '0.86', '0.0446', 4)
    ,@('011499', '富国沪深300基本面精选股票型证券投资基金C', '0.92', '90.96', '3.96', '0.0364', 7)
    ,@('005357', '富国国企改革灵活配置混合', '1.13', '87.21', '2.78', '0.0314', 2)
    ,@('002271', '招商安弘灵活配置混合', '0.50', '72.34', '5.16', '0.0258', 2)
    ,@('009908', '湘财长泽灵活配置混合C', '0.46', '80.11', '4.94', '0.0227', 3)
    ,@('519025', '海富通领先成长混合', '1.15', '93.44', '1.95', '0.0224', 10)
    ,@('002505', '鹏华金鼎灵活配置混合C', '0.26', '77.53', '6.72', '0.0175', 3)
    ,@('005672', '新疆前海联合研究优选灵活配置混合C', '0.21', '77.46', '3.65', '0.0077', 8)
    ,@('011149', '创金合信ESG责任投资股票A', '0.16', '87.53', '4.23', '0.0068', 4)
    ,@('001791', '大成绝对收益策略混合A', '0.31', '62.72', '2.07', '0.0064', 7)
    ,@('011150', '创金合信ESG责任投资股票C', '0.08', '87.53', '4.23', '0.0034', 4)
    ,@('011161', '富国质量成长6个月持有期混合C', '0.12', '91.55', '2.12', '0.0025', 10)
    ,@('163821', '中银沪深300等权重指数(LOF)', '0.46', '91.25', '0.48', '0.0022', 8)
    ,@('001792', '大成绝对收益策略混合C', '0.03', '62.72', '2.07', '0.0006', 7)
)

for ($i = 0; $i -lt $fundData.Count; $i++) {
    $r = $i + 2
    $row = $fundData[$i]
    $newSheet.Cells.Item($r, 1).Value2 = $i
    $newSheet.Cells.Item($r, 2).Value2 = $row[0]
    $newSheet.Cells.Item($r, 3).Value2 = $row[1]
    $newSheet.Cells.Item($r, 4).Value2 = $row[2]
    $newSheet.Cells.Item($r, 5).Value2 = $row[3]
    $newSheet.Cells.Item($r, 6).Value2 = $row[4]
    $newSheet.Cells.Item($r, 7).Value2 = $row[5]
    $newSheet.Cells.Item($r, 8).Value2 = $row[6]
}

# ---- Step 3: insert the 2022-Q1 summary row into "总计" ----
$totalSheet.Rows.Item(2).Insert()
# clone the (now shifted-down) old row 2 formatting onto the fresh row 2
$totalSheet.Range("A3:D3").Copy($totalSheet.Range("A2:D2"))

# bump the existing index column (A) of the old rows, now at rows 3..7, by +1
for ($r = 7; $r -ge 3; $r--) {
    $oldIndex = $totalSheet.Cells.Item($r, 1).Value2
    $totalSheet.Cells.Item($r, 1).Value2 = $oldIndex + 1
}

$totalSheet.Cells.Item(2, 1).Value2 = 0
$totalSheet.Cells.Item(2, 2).Value2 = "2022-Q1"
$totalSheet.Cells.Item(2, 3).Value2 = 34
$totalSheet.Cells.Item(2, 4).Value2 = 11.14

Write-Host "edit.ps1 completed"
